# Updated: st 23. 12. 2021
# Apply corrected AgTests (col F) / AgPosit (col G) figures to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(406, 6).Value = 171775
$ws.Cells.Item(421, 6).Value = 153619
$ws.Cells.Item(424, 6).Value = 266819
$ws.Cells.Item(426, 6).Value = 107850
$ws.Cells.Item(433, 6).Value = 87472
$ws.Cells.Item(440, 6).Value = 73950
$ws.Cells.Item(447, 6).Value = 67176
$ws.Cells.Item(478, 6).Value = 55275
$ws.Cells.Item(484, 6).Value = 8350
$ws.Cells.Item(506, 6).Value = 11113
$ws.Cells.Item(559, 6).Value = 22108
$ws.Cells.Item(586, 6).Value = 33893
$ws.Cells.Item(624, 6).Value = 51121
$ws.Cells.Item(626, 6).Value = 20051
$ws.Cells.Item(627, 6).Value = 33672
$ws.Cells.Item(630, 6).Value = 46508
$ws.Cells.Item(630, 7).Value = 2950
$ws.Cells.Item(631, 6).Value = 41800
$ws.Cells.Item(631, 7).Value = 2747
$ws.Cells.Item(632, 6).Value = 44071
$ws.Cells.Item(632, 7).Value = 2647
$ws.Cells.Item(633, 6).Value = 24011
$ws.Cells.Item(633, 7).Value = 1933
$ws.Cells.Item(634, 6).Value = 46497
$ws.Cells.Item(634, 7).Value = 2179
$ws.Cells.Item(635, 6).Value = 82994
$ws.Cells.Item(635, 7).Value = 3693
$ws.Cells.Item(636, 6).Value = 49779
$ws.Cells.Item(636, 7).Value = 2335
$ws.Cells.Item(637, 6).Value = 43342
$ws.Cells.Item(637, 7).Value = 2090
$ws.Cells.Item(638, 6).Value = 37353
$ws.Cells.Item(638, 7).Value = 1952
$ws.Cells.Item(639, 6).Value = 40398
$ws.Cells.Item(639, 7).Value = 1958
$ws.Cells.Item(640, 6).Value = 19470
$ws.Cells.Item(640, 7).Value = 1214
$ws.Cells.Item(641, 6).Value = 33601
$ws.Cells.Item(641, 7).Value = 1368
$ws.Cells.Item(642, 6).Value = 67219
$ws.Cells.Item(642, 7).Value = 2379
$ws.Cells.Item(643, 6).Value = 43150
$ws.Cells.Item(643, 7).Value = 1649
$ws.Cells.Item(644, 6).Value = 36235
$ws.Cells.Item(644, 7).Value = 1475
$ws.Cells.Item(645, 6).Value = 35168
$ws.Cells.Item(645, 7).Value = 1294
$ws.Cells.Item(646, 6).Value = 35743
$ws.Cells.Item(647, 6).Value = 16106
$ws.Cells.Item(647, 7).Value = 903
$ws.Cells.Item(648, 6).Value = 29922
$ws.Cells.Item(648, 7).Value = 1060
$ws.Cells.Item(649, 6).Value = 61900
$ws.Cells.Item(649, 7).Value = 1796
$ws.Cells.Item(650, 6).Value = 37550
$ws.Cells.Item(650, 7).Value = 1171
$ws.Cells.Item(651, 6).Value = 34686
$ws.Cells.Item(651, 7).Value = 1022
$ws.Cells.Item(652, 6).Value = 34277
$ws.Cells.Item(652, 7).Value = 1052
$ws.Cells.Item(653, 6).Value = 33064
$ws.Cells.Item(653, 7).Value = 978
$ws.Cells.Item(654, 6).Value = 13688
$ws.Cells.Item(654, 7).Value = 664
$ws.Cells.Item(655, 6).Value = 23872
$ws.Cells.Item(655, 7).Value = 769
$ws.Cells.Item(656, 6).Value = 48556
$ws.Cells.Item(656, 7).Value = 1155
$ws.Cells.Item(657, 6).Value = 31807
$ws.Cells.Item(657, 7).Value = 768
$ws.Cells.Item(658, 6).Value = 18846
$ws.Cells.Item(658, 7).Value = 466
